# Fill in X5/Y5 (which were missing) and append rows 6-11 of new random-walk
# prediction data, matching the structure of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Complete row 5, which was missing the last two columns ----
$ws.Range("X5").Value = 0.18999999999999773
$ws.Range("Y5").Value = "Up"

# ---- Data for the newly appended rows 6-11 ----
# Columns: A Date(serial) B ScoreFinal C Verdict D..O zeros P Method
#          Q RSI R PEG S 200Moving% T 50Moving% U PriceBook V Dividend
#          W Bollinger X PriceChange Y UpDown
$rows = @(
    @{ Row=6;  Date=42650.338379629633; Score=-7;  Verdict="Sell";        RSI=38.48959524716075; S=0.1046; T=0.0345; X=0.18999999999999773; Y="Up" },
    @{ Row=7;  Date=42650.339618055557; Score=6;   Verdict="Neutral";     RSI=38.48959524716075; S=0.1046; T=0.0345; X=0.18999999999999773; Y="Up" },
    @{ Row=8;  Date=42650.34878472222;  Score=-18; Verdict="Strong Sell"; RSI=38.48959524716075; S=0.1046; T=0.0345; X=0.18999999999999773; Y="Up" },
    @{ Row=9;  Date=42650.359050925923; Score=12;  Verdict="Buy";         RSI=38.48959524716075; S=0.1046; T=0.0345; X=0.18999999999999773; Y="Up" },
    @{ Row=10; Date=42650.361481481479; Score=-13; Verdict="Sell";        RSI=38.48959524716075; S=0.1046; T=0.0345; X=0.18999999999999773; Y="Up" },
    @{ Row=11; Date=42650.36310185185;  Score=12;  Verdict="Buy";         RSI=37.799019424898844; S=0.1046; T=0.0343; X=$null; Y=$null }
)

foreach ($r in $rows) {
    $row = $r.Row

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = $r.Date
    $cellA.NumberFormat = "m/d/yy h:mm"

    $ws.Cells.Item($row, 2).Value = $r.Score
    $ws.Cells.Item($row, 3).Value = $r.Verdict

    for ($col = 4; $col -le 15; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }

    $ws.Cells.Item($row, 16).Value = "Random"
    $ws.Cells.Item($row, 17).Value = $r.RSI
    $ws.Cells.Item($row, 18).Value = 0

    $cellS = $ws.Cells.Item($row, 19)
    $cellS.Value = $r.S
    $cellS.NumberFormat = "0.00%"

    $cellT = $ws.Cells.Item($row, 20)
    $cellT.Value = $r.T
    $cellT.NumberFormat = "0.00%"

    $ws.Cells.Item($row, 21).Value = 4.82
    $ws.Cells.Item($row, 22).Value = 2.2799999999999998
    $ws.Cells.Item($row, 23).Value = 0

    if ($null -ne $r.X) {
        $ws.Cells.Item($row, 24).Value = $r.X
    }
    if ($null -ne $r.Y) {
        $ws.Cells.Item($row, 25).Value = $r.Y
    }
}

# ---- Match the selection recorded in the diff (active cell B7) ----
$ws.Range("B7").Select()
